$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Include from hp.owl")

$ws.Range("A2").Value = "HP:0003577"
$ws.Range("B2").Value = "Congenital onset"

$ws.Range("A3").Value = "HP:0030674"
$ws.Range("B3").Value = "Antenatal onset"

$ws.Range("A4").Value = "HP:0011460"
$ws.Range("B4").Value = "Embryonal onset"

$ws.Range("A5").Value = "HP:0011461"
$ws.Range("B5").Value = "Fetal onset"

$ws.Range("A6").Value = "HP:0003623"
$ws.Range("B6").Value = "Neonatal onset"

$ws.Range("A7").Value = "HP:0003593"
$ws.Range("B7").Value = "Infantile onset"

$ws.Range("A8").Value = "HP:0011463"
$ws.Range("B8").Value = "Childhood onset"

$ws.Range("A9").Value = "HP:0003621"
$ws.Range("B9").Value = "Juvenile onset"

$ws.Range("A10").Value = "HP:0410280"
$ws.Range("B10").Value = "Pediatric onset"

$ws.Range("A11").Value = "HP:0003581"
$ws.Range("B11").Value = "Adult onset"

$ws.Range("A12").Value = "HP:0011462"
$ws.Range("B12").Value = "Young adult onset"

$ws.Range("A13").Value = "HP:0003596"
$ws.Range("B13").Value = "Middle age onset"

$ws.Range("A14").Value = "HP:0003584"
$ws.Range("B14").Value = "Late onset"
